$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "npc_dist_lt_dut_start" column (F)
$ws.Columns("F").Insert()

# Populate the new column: value row first, then header row, so that the
# shared-string table gets "[5..30]kph" before "npc_speed_increase"
$ws.Range("F3").Value = "[5..30]kph"
$ws.Range("F2").Value = "npc_speed_increase"

# Match the column width used elsewhere in the sheet
$ws.Columns("F").ColumnWidth = 21.5

# Update the selected cell
$ws.Range("F2").Select() | Out-Null
